# Auto-update predictions and index for 2025-10-22
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = "47/51 Win Tips"
$ws.Range("F3").Value = "'92"

$ws.Range("E4").Value = "38/44 Win Tips"
$ws.Range("F4").Value = "'86"

$ws.Range("E5").Value = "37/44 Win Tips"
$ws.Range("F5").Value = "'84"

$ws.Range("E6").Value = "35/45 Win Tips"
$ws.Range("F6").Value = "'78"

$ws.Range("E8").Value = "31/39 Win Tips"
$ws.Range("F8").Value = "'79"

$ws.Range("E10").Value = "22/49 Win Tips"

$ws.Range("E11").Value = "21/35 Win Tips"
$ws.Range("F11").Value = "'60"
